$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F3").Value = 1342
$ws.Range("F4").Value = 1148
$ws.Range("F5").Value = 1038
$ws.Range("F6").Value = 1822
$ws.Range("F7").Value = 577
$ws.Range("F8").Value = 1216
$ws.Range("F11").Value = 128
$ws.Range("F12").Value = 306
$ws.Range("F13").Value = 83
$ws.Range("F14").Value = 91
$ws.Range("F15").Value = 711
$ws.Range("F16").Value = 186
$ws.Range("F17").Value = 108
$ws.Range("F18").Value = 29
$ws.Range("F20").Value = 332
$ws.Range("F21").Value = 167
$ws.Range("F25").Value = 165
$ws.Range("F28").Value = 322
$ws.Range("F30").Value = 48
$ws.Range("F31").Value = 282
$ws.Range("F34").Value = 407

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F4").Value = 323
$ws.Range("F7").Value = 260
$ws.Range("F12").Value = 22

$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F2").Value = 314

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value = 314
$ws.Range("F4").Value = 1342
$ws.Range("F5").Value = 1148
$ws.Range("F6").Value = 1038
$ws.Range("F7").Value = 1822
$ws.Range("F8").Value = 577
$ws.Range("F9").Value = 1216
$ws.Range("F13").Value = 128
$ws.Range("F14").Value = 306
$ws.Range("F15").Value = 83
$ws.Range("F16").Value = 91
$ws.Range("F17").Value = 711
$ws.Range("F18").Value = 186
$ws.Range("F19").Value = 108
$ws.Range("F21").Value = 29
$ws.Range("F22").Value = 323
$ws.Range("F25").Value = 332
$ws.Range("F27").Value = 260
$ws.Range("F28").Value = 260
$ws.Range("F29").Value = 167
$ws.Range("F33").Value = 165
$ws.Range("F36").Value = 322
$ws.Range("F40").Value = 48
$ws.Range("F41").Value = 282
$ws.Range("F47").Value = 22
$ws.Range("F48").Value = 407

